# [FIX] replaced doublequotes for single
#
# The "Cluster" sheet holds rule/formula definitions (JS-like snippets such as
# ["GROUP_5","GROUP_6", ...].indexOf(model.countryGroup) >= 0) in column D.
# These snippets used double quotes around the string literals; the author
# changed them to single quotes. Reproduce that with a global find/replace
# across the used range, exactly like using Excel's Find & Replace (Ctrl+H)
# with " -> ' would do.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace every double quote with a single quote across all used cells.
$ws.UsedRange.Replace('"', "'") | Out-Null

# The author's last selection/scroll position ended up on D82 (row 64 visible
# at the top) after making this edit - reflect that in the saved view state.
$ws.Range("D82").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
